{"js": "// The edit updates the worksheet's date line and every three-digit x\n// one-digit multiplication \"problem=answer\" string in the practice table\n// with a new problem/answer pair. Every old string in the mapping below is\n// unique within the document, so a straight search-and-replace per pair is\n// safe and unambiguous.\nconst replacements = [\n  [\"2024-03-18 Monday\", \"2024-03-19 Tuesday\"],\n  [\"719\u00d72=1438\", \"946\u00d78=7568\"],\n  [\"351\u00d77=2457\", \"845\u00d74=3380\"],\n  [\"395\u00d78=3160\", \"706\u00d79=6354\"],\n  [\"220\u00d74=880\", \"955\u00d75=4775\"],\n  [\"387\u00d76=2322\", \"279\u00d78=2232\"],\n  [\"999\u00d78=7992\", \"253\u00d78=2024\"],\n  [\"617\u00d78=4936\", \"367\u00d74=1468\"],\n  [\"830\u00d73=2490\", \"323\u00d75=1615\"],\n  [\"660\u00d75=3300\", \"935\u00d78=7480\"],\n  [\"662\u00d77=4634\", \"801\u00d79=7209\"],\n  [\"141\u00d72=282\", \"936\u00d75=4680\"],\n  [\"330\u00d73=990\", \"692\u00d73=2076\"],\n  [\"969\u00d73=2907\", \"892\u00d77=6244\"],\n  [\"985\u00d73=2955\", \"241\u00d78=1928\"],\n  [\"223\u00d74=892\", \"977\u00d74=3908\"],\n  [\"361\u00d78=2888\", \"446\u00d78=3568\"],\n  [\"604\u00d78=4832\", \"217\u00d72=434\"],\n  [\"216\u00d76=1296\", \"332\u00d73=996\"],\n  [\"794\u00d78=6352\", \"370\u00d78=2960\"],\n  [\"268\u00d76=1608\", \"376\u00d74=1504\"],\n  [\"168\u00d72=336\", \"422\u00d79=3798\"],\n  [\"259\u00d73=777\", \"931\u00d76=5586\"],\n  [\"353\u00d78=2824\", \"150\u00d79=1350\"],\n  [\"393\u00d77=2751\", \"493\u00d73=1479\"],\n  [\"592\u00d74=2368\", \"839\u00d77=5873\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The edit updates the worksheet's date line and every three-digit x\n# one-digit multiplication \"problem=answer\" string in the practice table\n# with a new problem/answer pair. Every old string in the mapping below is\n# unique within the document, so a straight Find/Replace per pair is safe\n# and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-03-18 Monday\", \"2024-03-19 Tuesday\"),\n  @(\"719\u00d72=1438\", \"946\u00d78=7568\"),\n  @(\"351\u00d77=2457\", \"845\u00d74=3380\"),\n  @(\"395\u00d78=3160\", \"706\u00d79=6354\"),\n  @(\"220\u00d74=880\", \"955\u00d75=4775\"),\n  @(\"387\u00d76=2322\", \"279\u00d78=2232\"),\n  @(\"999\u00d78=7992\", \"253\u00d78=2024\"),\n  @(\"617\u00d78=4936\", \"367\u00d74=1468\"),\n  @(\"830\u00d73=2490\", \"323\u00d75=1615\"),\n  @(\"660\u00d75=3300\", \"935\u00d78=7480\"),\n  @(\"662\u00d77=4634\", \"801\u00d79=7209\"),\n  @(\"141\u00d72=282\", \"936\u00d75=4680\"),\n  @(\"330\u00d73=990\", \"692\u00d73=2076\"),\n  @(\"969\u00d73=2907\", \"892\u00d77=6244\"),\n  @(\"985\u00d73=2955\", \"241\u00d78=1928\"),\n  @(\"223\u00d74=892\", \"977\u00d74=3908\"),\n  @(\"361\u00d78=2888\", \"446\u00d78=3568\"),\n  @(\"604\u00d78=4832\", \"217\u00d72=434\"),\n  @(\"216\u00d76=1296\", \"332\u00d73=996\"),\n  @(\"794\u00d78=6352\", \"370\u00d78=2960\"),\n  @(\"268\u00d76=1608\", \"376\u00d74=1504\"),\n  @(\"168\u00d72=336\", \"422\u00d79=3798\"),\n  @(\"259\u00d73=777\", \"931\u00d76=5586\"),\n  @(\"353\u00d78=2824\", \"150\u00d79=1350\"),\n  @(\"393\u00d77=2751\", \"493\u00d73=1479\"),\n  @(\"592\u00d74=2368\", \"839\u00d77=5873\")\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format, ReplaceWith,\n#   Replace) -- Wrap:=1 is wdFindContinue, Replace:=2 is wdReplaceAll.\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
